$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45978, 8784.3061304421, 7422.16946448065, 12603.4, 4362.64253292383, -34.1078334414802),
    @(45979, 9584.56161281496, 8840.28375401519, 9003.4, 5435.88294533162, 219.698612472783),
    @(45980, 10327.6048500273, 9844.86614253143, 9003.4, 5963.05527828457, 283.521725867334),
    @(45981, 10491.7680288479, 9895.38760223554, 9003.4, 6075.5153766779, 290.312624121393),
    @(45982, 10534.8266030694, 9560.92987631162, 9003.4, 6113.06087672843, 277.941281376669),
    @(45983, 4000.30754528634, 7192.60797284348, 9003.4, 5477.70624761572, 152.788092519134),
    @(45984, 3853.28814954241, 7164.69858329255, 9003.4, 5469.70274819967, 151.291722145509),
    @(45985, 10980.8696790314, 10647.3889751765, 9003.4, 6409.18310587767, 335.548836710592),
    @(45986, 10980.8696790314, 10953.6284220058, 9003.4, 6409.18310587767, 348.308813661812),
    @(45987, 10980.8696790314, 11080.0905991367, 9003.4, 6409.18310587767, 353.578071042266),
    @(45988, 10980.8696790314, 10698.6529613641, 9003.4, 6409.18310587767, 337.684836135075),
    @(45989, 10980.8696790314, 9684.13321846537, 9003.4, 6409.18310587767, 295.41318018096),
    @(45990, 4442.42485201967, 6618.19460851034, 9003.4, 6032.96273272966, 151.989889218333),
    @(45991, 4286.74551609542, 6382.33152571853, 9003.4, 6024.47415076784, 141.808569853599)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
}
